# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a failed
# handback transform for the c22f96d5-... file:
#   - The shared "Ready for handoff" status text becomes
#     "Handback transform failed" (affects Overview!E3/F3 and the
#     Status column on both the zh-cn and de-de sheets, since they all
#     point at the same status string).
#   - The Error Detail column (P) on the zh-cn and de-de sheets is filled
#     in with an explanatory message, and that column is widened to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status text. Overview!E3 and F3 (zh-cn/de-de status for the
# c22f96d5 row) currently read "Ready for handoff" - change to reflect the
# failed handback transform.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# The Status column (C) on the zh-cn / de-de detail sheets shares the same
# text for row 3 (the c22f96d5 file).
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Fill in the Error Detail (column P) for row 3 on each locale sheet.
$zhcn.Range("P3").Value = "Handback file name: bbaogo53.pfh is different with handoff file name: c22f96d5-f6af-4fc1-a4f1-ff60ec745def.8f750ec828fe68a9ea64b2c8e38c20a6acf0854f.zh-cn."
$dede.Range("P3").Value = "Handback file name: bbaogo53.pfh is different with handoff file name: c22f96d5-f6af-4fc1-a4f1-ff60ec745def.8f750ec828fe68a9ea64b2c8e38c20a6acf0854f.de-de."

# Widen column P (Error Detail) on both sheets to width 40 so the message
# is readable (raw OOXML column width = ColumnWidth + 5/6).
$zhcn.Columns.Item(16).ColumnWidth = 40 - 0.8333333333333334
$dede.Columns.Item(16).ColumnWidth = 40 - 0.8333333333333334
